$wb = $excel.ActiveWorkbook

# --- Metadata sheet ---
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B2").Value = "http://linuxforhealth.org/fhir/cdm/ValueSet/patient-status"
$meta.Range("B3").Value = "8.0.0"
$meta.Range("B8").Value = "2022-11-10T16:00:46+00:00"
$meta.Range("B9").Value = "LinuxForHealth Team"

# --- Include from Patient Status C sheet ---
$codeSheet = $wb.Worksheets.Item("Include from Patient Status C")
$codeSheet.Range("B4").Value = "http://linuxforhealth.org/fhir/cdm/CodeSystem/patient-status"
